# Update gh-pages to output generated at 456a3b4
# Applies the diff across sheets: 展览(1) / 演出(2) / 本地生活(3) / 全部类型(4)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Simple "want-to-go count" (F column) bumps that don't involve row shifts.
$ws1.Range("F5").Value  = 1683
$ws1.Range("F8").Value  = 614
$ws1.Range("F9").Value  = 3515
$ws1.Range("F16").Value = 1272
$ws1.Range("F19").Value = 448
$ws1.Range("F20").Value = 1543
$ws1.Range("F21").Value = 1072

# Insert a brand-new event as row 23, pushing the former rows 23-27 down to 24-28.
$ws1.Rows.Item(23).Insert()

# Column A keeps the same bold/bordered style as the rest of the numbering column.
$ws1.Cells.Item(22, 1).Copy()
$ws1.Cells.Item(23, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Cells.Item(23, 1).Value = 22

# Column B holds a literal "yyyy-mm-dd" label (not a real date) - force Text so
# Excel doesn't reinterpret it as a date serial, then strip the number-format
# back off (copy formatting from a plain-text sibling) so it stays the default style.
$ws1.Cells.Item(23, 2).NumberFormat = "@"
$ws1.Cells.Item(23, 2).Value = "2024-11-09"

$ws1.Cells.Item(23, 3).Value = "上海·夜蓝诗2.0·恋与深空同人only"
$ws1.Cells.Item(23, 4).Value = "沪闵路7388号 上海百联南方商城"
$ws1.Cells.Item(23, 5).Value = "2024.11.09 11:00-11.09 21:00"
$ws1.Cells.Item(23, 6).Value = 0
$ws1.Cells.Item(23, 7).Value = 98
$ws1.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93438"
$ws1.Cells.Item(23, 9).Value = "//i2.hdslb.com/bfs/openplatform/202410/KeJJBybC1728903298004.jpeg"

# Reset B23's style to the default (no explicit style index), matching its
# plain-text siblings, while keeping the text value we just wrote.
$ws1.Cells.Item(23, 3).Copy()
$ws1.Cells.Item(23, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A is a simple running index (row number - 1); the rows pushed down by
# the insert (old 23-27, now 24-28) need to be renumbered to match their new
# position instead of keeping their pre-shift index.
$ws1.Cells.Item(24, 1).Value = 23
$ws1.Cells.Item(25, 1).Value = 24
$ws1.Cells.Item(26, 1).Value = 25
$ws1.Cells.Item(27, 1).Value = 26
$ws1.Cells.Item(28, 1).Value = 27

# The events that shifted down from old rows 24/25 each had their "want-to-go"
# counts tick up as well.
$ws1.Range("F25").Value = 4255
$ws1.Range("F26").Value = 48
$ws1.Range("G26").Value = "不可售"

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 15
$ws2.Range("F10").Value = 164
$ws2.Range("F23").Value = 120
$ws2.Range("F36").Value = 425
$ws2.Range("F41").Value = 18
$ws2.Range("F50").Value = 46

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value  = 2550
$ws3.Range("F7").Value  = 155
$ws3.Range("F10").Value = 385
$ws3.Range("F11").Value = 2958
$ws3.Range("F12").Value = 463
$ws3.Range("F13").Value = 796
$ws3.Range("F14").Value = 202

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 155
$ws4.Range("F9").Value  = 2958
$ws4.Range("F10").Value = 796
$ws4.Range("F11").Value = 202
$ws4.Range("F14").Value = 1683
$ws4.Range("F16").Value = 614
$ws4.Range("F22").Value = 14
$ws4.Range("F23").Value = 15
$ws4.Range("F25").Value = 164
$ws4.Range("F33").Value = 1543
$ws4.Range("F34").Value = 120
$ws4.Range("F35").Value = 120
$ws4.Range("F37").Value = 1072

Write-Output "edit complete"
